$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  4/28/2025  Through  5/4/2025"

# --- Row 14: 28-Day %Chg for Murder ---
$ws.Range("N14").Value = -60

# --- Row 15 (Rape): restructure - C/D/E become blank placeholders, F/G/H shift ---
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 400
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 33.333333333333
$ws.Range("N15").Value = 50

# --- Row 16 updates ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -38.095238095238
$ws.Range("I16").Value = 79
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = -16.842105263157
$ws.Range("L16").Value = 1.282051282051
$ws.Range("M16").Value = -13.186813186813
$ws.Range("N16").Value = -69.379844961240

# --- Row 17 updates ---
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 8.333333333333
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 19.354838709677
$ws.Range("I17").Value = 145
$ws.Range("J17").Value = 143
$ws.Range("K17").Value = 1.398601398601
$ws.Range("L17").Value = 9.848484848484
$ws.Range("M17").Value = 66.666666666666
$ws.Range("N17").Value = 64.772727272727

# --- Row 18 updates ---
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -42.857142857142
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -26.666666666666
$ws.Range("I18").Value = 89
$ws.Range("J18").Value = 64
$ws.Range("K18").Value = 39.0625
$ws.Range("L18").Value = 8.536585365853
$ws.Range("M18").Value = -14.423076923076
$ws.Range("N18").Value = -83.876811594202

# --- Row 19 updates ---
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -19.047619047619
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -8.474576271186
$ws.Range("I19").Value = 254
$ws.Range("J19").Value = 324
$ws.Range("K19").Value = -21.604938271604
$ws.Range("L19").Value = 22.705314009661
$ws.Range("M19").Value = 124.778761061947
$ws.Range("N19").Value = 34.391534391534

# --- Row 20 updates ---
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = -18.181818181818
$ws.Range("I20").Value = 164
$ws.Range("J20").Value = 154
$ws.Range("K20").Value = 6.493506493506
$ws.Range("L20").Value = 1.863354037267
$ws.Range("M20").Value = 110.25641025641
$ws.Range("N20").Value = -74.050632911392

# --- Row 21 updates ---
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = -13.461538461538
$ws.Range("F21").Value = 152
$ws.Range("G21").Value = 157
$ws.Range("H21").Value = -3.184713375796
$ws.Range("I21").Value = 745
$ws.Range("J21").Value = 794
$ws.Range("K21").Value = -6.171284634760
$ws.Range("L21").Value = 11.194029850746
$ws.Range("M21").Value = 53.925619834710
$ws.Range("N21").Value = -56.986143187067

# --- Row 22 updates ---
$ws.Range("G22").Value = 1
$ws.Range("M22").Value = -37.5

# --- Row 23 updates ---
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -75
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -63.636363636363
$ws.Range("I23").Value = 34
$ws.Range("J23").Value = 43
$ws.Range("K23").Value = -20.930232558139
$ws.Range("L23").Value = -29.166666666666
$ws.Range("M23").Value = 61.904761904761

# --- Row 24 updates ---
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -12.903225806451
$ws.Range("F24").Value = 132
$ws.Range("G24").Value = 90
$ws.Range("H24").Value = 46.666666666666
$ws.Range("I24").Value = 495
$ws.Range("J24").Value = 510
$ws.Range("K24").Value = -2.941176470588
$ws.Range("L24").Value = -4.807692307692
$ws.Range("M24").Value = 72.473867595818

# --- Row 25 updates ---
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -33.333333333333
$ws.Range("I25").Value = 171
$ws.Range("J25").Value = 220
$ws.Range("K25").Value = -22.272727272727
$ws.Range("L25").Value = -21.559633027522

# --- Row 26 updates ---
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = 42.857142857142
$ws.Range("I26").Value = 220
$ws.Range("J26").Value = 168
$ws.Range("K26").Value = 30.952380952381
$ws.Range("L26").Value = 17.021276595744
$ws.Range("M26").Value = 1.851851851851

# --- Row 27 (UCR Rape*): restructure - C becomes blank, D/E get real values, F/G/H/J/K/L updated ---
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 200
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 15.384615384615

# --- Row 28 (Other Sex Crimes) updates ---
$ws.Range("C28").Value = 4
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 8
$ws.Range("H28").Value = 14.285714285714
$ws.Range("I28").Value = 26
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = 13.043478260869
$ws.Range("L28").Value = 52.941176470588

# --- Row 29 (Shooting Vic.): D/E become blank placeholders ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"

# --- Row 30 (Shooting Inc.): D/E become blank placeholders ---
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
